# WS_holdings.xlsx update:
#  - Roll the "as of" date in the confidential disclaimer from 2021-05-19 to 2021-05-20
#  - Refresh the Weight (D) and Percent Change (E) figures for rows 2-13
#
# The sheet ships protected (sheetProtection ... sheet="1" objects="1" scenarios="1"
# formatColumns="0" formatRows="0"), so it has to be unprotected before the cells can
# be written, and re-protected afterwards to leave the sheet in the same locked state
# it started in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Disclaimer text (A16) - bump the model-holdings-as-of date by one day.
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-20 for illustrative purposes only and are subject to change."

# Weight (D) / Percent Change (E) columns, rows 2-13
$ws.Range("D2").Value = 0.03103561053594975
$ws.Range("E2").Value = 0.01138716356107627

$ws.Range("D3").Value = 0.02385483070971028
$ws.Range("E3").Value = 0.006638217164532945

$ws.Range("D4").Value = 0.05265375756360795
$ws.Range("E4").Value = 0.007226107226107281

$ws.Range("D5").Value = 0.1389552372066338
$ws.Range("E5").Value = 0.01309328968903434

$ws.Range("D6").Value = 0.03179476899033249
$ws.Range("E6").Value = -0.0007077140835102247

$ws.Range("D7").Value = 0.1167838154011579
$ws.Range("E7").Value = 0.007851347814707976

$ws.Range("D8").Value = 0.1023038515770016
$ws.Range("E8").Value = 0.004088459394164579

$ws.Range("D9").Value = 0.02964688843703846
$ws.Range("E9").Value = 0.0006209894431794538

$ws.Range("D10").Value = 0.1277417364356831
$ws.Range("E10").Value = 0.001130369253956376

$ws.Range("D11").Value = 0.24330372667504
$ws.Range("E11").Value = 0.02020296936666033

$ws.Range("D12").Value = 0.1019257764678447
$ws.Range("E12").Value = 0.01864339547798499

$ws.Range("E13").Value = 0.01100280391667519

# Restore the protected state.
$ws.Protect()
